# Correct the ZEV biosample "treatment" column: cells that were labeled
# "EtOH" should instead read "mockEstradiol" (the "Estradiol" treatment
# values are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)   # column J = treatment
    $val = $cell.Value()
    if ($val -eq "EtOH") {
        $cell.Value = "mockEstradiol"
    }
}

# Update the sheet's active selection to match the author's last cursor
# position when they saved the workbook.
$null = $ws.Range("O12").Select()
